$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B5").Value = 0
$ws.Range("C5").Value = "Will McLain"
$ws.Range("D5").Value = "Covered the working agreement for the propsal"
